$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: the shared string used by B10/C10 changed text (long Objectives
#     paragraph replaced by the responsible-professor line). ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("C10").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 60

# --- Rows 13-25 are restructured/reshuffled substantially; delete them all
#     and rebuild rows 13-23 fresh with the final content. ---
$ws.Range("A13:A25").EntireRow.Delete()

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2019"
$ws.Range("C13").Value = "01/01/2019"
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Characteristics of biological material; Kinetics of fermentative processes; Operation modes of fermentative processes; Sterilization in bioprocess."
$ws.Range("C14").Value = "Characteristics of biological material; Kinetics of fermentative processes; Operation modes of fermentative processes; Sterilization in bioprocess."
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("C15").Value = "1112574 - Inês Conceição Roberto"
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Characteristics of biological material - Types of microorganisms, cell structure and morphology, nutrition and microbial growth. Kinetics of fermentative processes - Evaluation of kinetic profile of batch cultures, definition and calculation of fermentative parameters (rates and yields), kinetic models for cell growth (MONOD) and product formation (GADEN). Operation modes of fermentative processes. Major configurations of bioreactors, characteristics and mathematical equations for batch, fed-batch and continuous operations. Sterilization in fermentation process – general aspects on sterilization and disinfection in bioprocess, methods for medium and air sterilization, kinetics of thermal death of microorganisms, calculation of sterilization time for batch and continuous process."
$ws.Range("C16").Value = "Characteristics of biological material - Types of microorganisms, cell structure and morphology, nutrition and microbial growth. Kinetics of fermentative processes - Evaluation of kinetic profile of batch cultures, definition and calculation of fermentative parameters (rates and yields), kinetic models for cell growth (MONOD) and product formation (GADEN). Operation modes of fermentative processes. Major configurations of bioreactors, characteristics and mathematical equations for batch, fed-batch and continuous operations. Sterilization in fermentation process – general aspects on sterilization and disinfection in bioprocess, methods for medium and air sterilization, kinetics of thermal death of microorganisms, calculation of sterilization time for batch and continuous process."
$ws.Range("A3").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("C18").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."
$ws.Range("C19").Value = "Os alunos serão avaliados formalmente por duas provas teóricas (P1 e P2). A ponderação das notas será de 50% para cada avaliação, ou seja: Média do período letivo normal = (P1+ P2)/2."
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Range("C20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("C21").Value = "Aos alunos que não obtiverem média igual ou maior que 5,0 será oferecido um programa de recuperação, que será avaliado por uma prova final (PF). Neste caso, a média final do aluno será: Média Final = (Média do período letivo normal + nota prova final) / 2. Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("A3").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Range("B23").Value = "LOT2004 -  Bioquímica  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOT2004 -  Bioquímica  (Requisito fraco)`n"
$ws.Range("B3").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 30

$excel.CutCopyMode = 0
